$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Featured Sessions")
$ws.Activate()

# Insert a new row above current row 3 (pushes existing rows 3-5 down to 4-6)
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = "Biometrical Journal Showcase - Editor's Selection"
$ws.Cells.Item(3, 2).Value = "Arne Bathke, Matthias Schmid"

$ws.Columns.Item(2).ColumnWidth = 24.3046875

$ws.Range("C9").Select()
